$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.784.03"
$ws.Range("E2").Value = "  +1.33%  "
$ws.Range("D3").Value = "1.663.77"
$ws.Range("E3").Value = "  +1.66%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "328.87"
$ws.Range("E5").Value = "  +7.74%  "
$ws.Range("E6").Value = "  -0.04%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3644"
$ws.Range("E7").Value = "  +0.77%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "47.26"
$ws.Range("E8").Value = "  +0.72%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3254"
$ws.Range("E9").Value = "  +0.90%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.135"
$ws.Range("E10").Value = "  +2.89%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07067"
$ws.Range("E11").Value = "  +2.52%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.001"
$ws.Range("E12").Value = "  -0.24%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.065"
$ws.Range("E13").Value = "  +2.86%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "19.46"
$ws.Range("E14").Value = "  +1.75%  "
$ws.Range("D15").Value = "1.667.12"
$ws.Range("E15").Value = "  +1.69%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.589"
$ws.Range("E16").Value = "  +1.17%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001047"
$ws.Range("E17").Value = "  +0.73%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06645"
$ws.Range("E18").Value = "  +2.08%  "
$ws.Range("E19").Value = "  -0.05%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "78.43"
$ws.Range("E20").Value = "  +2.56%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.926"
$ws.Range("E21").Value = "  +1.05%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "15.76"
$ws.Range("E22").Value = "  +0.26%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.50"
$ws.Range("E23").Value = "  +5.08%  "
$ws.Range("D24").Value = "24.815.35"
$ws.Range("E24").Value = "  +1.78%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.469"
$ws.Range("E25").Value = "  +2.43%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.435"
$ws.Range("E26").Value = "  +2.57%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "148.67"
$ws.Range("E27").Value = "  +3.44%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.67"
$ws.Range("E28").Value = "  +0.25%  "
$ws.Range("D29").Value = "1.847.37"
$ws.Range("E29").Value = "  +1.46%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "125.94"
$ws.Range("E30").Value = "  +1.10%  "
$ws.Range("E31").Value = "  +6.19%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.070"
$ws.Range("E32").Value = "  +0.15%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.671"
$ws.Range("E33").Value = "  +0.76%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08454"
$ws.Range("E34").Value = "  +0.76%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.645"
$ws.Range("E35").Value = "  -1.38%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "12.12"
$ws.Range("E36").Value = "  -1.50%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06236"
$ws.Range("E37").Value = "  +4.03%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.159"
$ws.Range("E38").Value = "  +1.01%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02277"
$ws.Range("E39").Value = "  +3.23%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.243"
$ws.Range("E40").Value = "  +3.70%  "
$ws.Range("E41").Value = "  +2.48%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.209"
$ws.Range("E42").Value = "  +0.99%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.001"
$ws.Range("E43").Value = "  -0.09%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5937"
$ws.Range("E44").Value = "  +1.76%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.41"
$ws.Range("E45").Value = "  +7.53%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.846"
$ws.Range("E46").Value = "  +3.51%  "
$ws.Range("E47").Value = "  +2.27%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "125.49"
$ws.Range("E48").Value = "  +3.46%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.948"
$ws.Range("E49").Value = "  +1.97%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06970"
$ws.Range("E50").Value = "  +1.02%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.194"
$ws.Range("E51").Value = "  +4.75%  "

Write-Host "Applied cryptos update"
